$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.011.34'
$ws.Range("E2").Value = '  -1.97%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.830.10'
$ws.Range("E3").Value = '  -1.11%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.50'
$ws.Range("E5").Value = '  -3.47%  '

$ws.Range("E6").Value = '  -0.04%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4639'
$ws.Range("E7").Value = '  -0.48%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3864'
$ws.Range("E8").Value = '  -1.05%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07839'
$ws.Range("E9").Value = '  -0.83%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9591'
$ws.Range("E10").Value = '  -2.30%  '

$ws.Range("E11").Value = '  -1.97%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.845.50'
$ws.Range("E12").Value = '  -0.47%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.676'
$ws.Range("E13").Value = '  -2.93%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.879'
$ws.Range("E14").Value = '  -1.84%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06864'
$ws.Range("E15").Value = '  -0.47%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '88.18'
$ws.Range("E16").Value = '  +0.64%  '

$ws.Range("E17").Value = '  -0.03%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009911'
$ws.Range("E18").Value = '  -1.30%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.61'
$ws.Range("E19").Value = '  -2.99%  '

$ws.Range("E20").Value = '  +0.04%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '28.018.83'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.286'
$ws.Range("E22").Value = '  -2.10%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.95'
$ws.Range("E23").Value = '  -3.28%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.087'
$ws.Range("E24").Value = '  -1.87%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.051.86'
$ws.Range("E25").Value = '  -0.73%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '154.91'
$ws.Range("E26").Value = '  +0.99%  '

$ws.Range("E27").Value = '  -1.75%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.650'
$ws.Range("E28").Value = '  -6.01%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.959'
$ws.Range("E29").Value = '  -3.23%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '118.30'
$ws.Range("E30").Value = '  +0.54%  '

$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.9350'
$ws.Range("E31").Value = '  -3.98%  '

$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09242'
$ws.Range("E32").Value = '  -1.79%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.246'
$ws.Range("E33").Value = '  -2.22%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.317'
$ws.Range("E34").Value = '  -2.32%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.309'
$ws.Range("E35").Value = '  -4.90%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05841'
$ws.Range("E36").Value = '  -5.12%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02120'
$ws.Range("E37").Value = '  -3.56%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.138'
$ws.Range("E38").Value = '  -1.88%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.740'
$ws.Range("E39").Value = '  +1.29%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5584'
$ws.Range("E40").Value = '  -2.49%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.857'
$ws.Range("E41").Value = '  -3.12%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1757'
$ws.Range("E42").Value = '  -2.47%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.07251'
$ws.Range("E43").Value = '  +1.50%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '11.56'
$ws.Range("E44").Value = '  -1.62%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5245'
$ws.Range("E45").Value = '  -2.75%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.142'
$ws.Range("E46").Value = '  -8.46%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.090'
$ws.Range("E47").Value = '  -11.43%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.820'
$ws.Range("E48").Value = '  -4.71%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '112.73'
$ws.Range("E49").Value = '  -2.85%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.000'
$ws.Range("E50").Value = '  -0.08%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.323'
$ws.Range("E51").Value = '  +0.45%  '
